$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.639
$ws.Range("AF5").Value = 0.833
$ws.Range("AF6").Value = 0.723
$ws.Range("AF7").Value = 0.785
$ws.Range("AF8").Value = 0.772
$ws.Range("AF9").Value = 0.667
$ws.Range("AF10").Value = 0.833
$ws.Range("AF11").Value = 0.833
$ws.Range("AF12").Value = 1.2
$ws.Range("AF13").Value = 1.667
